$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.615.81'
$ws.Range('E2').Value = '  -1.98%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.587.32'
$ws.Range('E3').Value = '  -2.24%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.17'
$ws.Range('E5').Value = '  -1.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.508'
$ws.Range('E6').Value = '  -2.65%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  -2.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0615'
$ws.Range('E9').Value = '  -1.85%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.52'
$ws.Range('E10').Value = '  -3.93%  '
$ws.Range('E11').Value = '  -1.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.810.13'
$ws.Range('E12').Value = '  -2.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.598.11'
$ws.Range('E13').Value = '  -1.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.02'
$ws.Range('E14').Value = '  -2.90%  '
$ws.Range('E15').Value = '  -3.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.86'
$ws.Range('E16').Value = '  +0.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.598.34'
$ws.Range('E17').Value = '  -1.99%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0726'
$ws.Range('E18').Value = '  -2.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '209.10'
$ws.Range('E19').Value = '  -3.06%  '
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.69'
$ws.Range('E21').Value = '  -3.26%  '
$ws.Range('E22').Value = '  -2.65%  '
$ws.Range('E23').Value = '  -2.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.84'
$ws.Range('E24').Value = '  -2.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.41'
$ws.Range('E25').Value = '  -1.91%  '
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.21'
$ws.Range('E27').Value = '  -0.91%  '
$ws.Range('E28').Value = '  -2.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.27'
$ws.Range('E29').Value = '  -2.17%  '
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('E31').Value = '  -1.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.22'
$ws.Range('E32').Value = '  -4.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.680'
$ws.Range('E33').Value = '  +22.66%  '
$ws.Range('E34').Value = '  -3.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.311.48'
$ws.Range('E35').Value = '  -2.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.45'
$ws.Range('E36').Value = '  -0.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.48'
$ws.Range('E37').Value = '  -5.24%  '
$ws.Range('E38').Value = '  -3.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.825'
$ws.Range('E39').Value = '  -4.05%  '
$ws.Range('E40').Value = '  +0.16%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.39'
$ws.Range('E41').Value = '  +2.89%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.791'
$ws.Range('E42').Value = '  -1.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.17'
$ws.Range('E43').Value = '  -2.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.54'
$ws.Range('E44').Value = '  -4.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.723.25'
$ws.Range('E45').Value = '  -2.02%  '
$ws.Range('E46').Value = '  -0.86%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.838'
$ws.Range('E48').Value = '  -9.68%  '
$ws.Range('E49').Value = '  -1.64%  '
$ws.Range('E50').Value = '  -2.31%  '
$ws.Range('E51').Value = '  -1.14%  '
